$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J6").Value = "foobar"
$ws.Range("K6").Value = "Psychrophrynella"
$ws.Range("L6").Value = "chirihampatu"

$ws.Range("J7").Select()
